$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) were reshuffled (dates + price columns permuted across rows),
# while all other columns (A,B,C,E,F,G,H,I,N,O,Q,R) stayed the same.
# Target values per row, taken from the diff:

$ws.Range("D2").Value = 44537
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1300
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = 1350
$ws.Range("P2").Value = 1350

$ws.Range("D3").Value = 44210
$ws.Range("J3").Value = 1450
$ws.Range("K3").Value = 1600
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1650
$ws.Range("P3").Value = 1650

$ws.Range("D4").Value = 44638
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2800
$ws.Range("M4").Value = 2650
$ws.Range("P4").Value = 2650

$ws.Range("D5").Value = 44200
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 1400
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1450
$ws.Range("P5").Value = 1450

$ws.Range("D6").Value = 44175
$ws.Range("J6").Value = 1400
$ws.Range("K6").Value = 1900
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1950
$ws.Range("P6").Value = 1950
